$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number-format on Price cells whose new value would otherwise be
# auto-parsed by Excel as a plain number (single decimal point, no thousand dots).
# This mirrors how Excel keeps a typed value as text (quote-prefix) instead of a number.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = '63.942.80'
$ws.Range("E2").Value = '  +1.48%  '
$ws.Range("D3").Value = '3.314.30'
$ws.Range("E3").Value = '  +6.03%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").Value = '600.42'
$ws.Range("E5").Value = '  +1.16%  '
$ws.Range("D6").Value = '143.16'
$ws.Range("E6").Value = '  +5.13%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").Value = '3.313.36'
$ws.Range("E8").Value = '  +6.37%  '
$ws.Range("E9").Value = '  +1.04%  '
$ws.Range("E10").Value = '  +2.77%  '
$ws.Range("D11").Value = '5.51'
$ws.Range("E11").Value = '  +5.39%  '
$ws.Range("D12").Value = '0.474'
$ws.Range("E12").Value = '  +3.22%  '
$ws.Range("E13").Value = '  +1.09%  '
$ws.Range("D14").Value = '34.87'
$ws.Range("E14").Value = '  +1.90%  '
$ws.Range("D15").Value = '3.859.30'
$ws.Range("E15").Value = '  +6.02%  '
$ws.Range("E16").Value = '  +1.13%  '
$ws.Range("D17").Value = '3.314.01'
$ws.Range("E17").Value = '  +6.12%  '
$ws.Range("D18").Value = '64.034.09'
$ws.Range("E18").Value = '  +1.69%  '
$ws.Range("E19").Value = '  +2.77%  '
$ws.Range("D20").Value = '481.96'
$ws.Range("E20").Value = '  +1.33%  '
$ws.Range("D21").Value = '14.28'
$ws.Range("E21").Value = '  +0.14%  '
$ws.Range("D22").Value = '0.743'
$ws.Range("E22").Value = '  +6.14%  '
$ws.Range("D23").Value = '8.01'
$ws.Range("E23").Value = '  +3.77%  '
$ws.Range("D24").Value = '13.65'
$ws.Range("E24").Value = '  +4.83%  '
$ws.Range("D25").Value = '84.76'
$ws.Range("E25").Value = '  -2.61%  '
$ws.Range("E26").Value = '  +0.17%  '
$ws.Range("D27").Value = '2.78'
$ws.Range("E27").Value = '  +2.14%  '
$ws.Range("E28").Value = '  +1.31%  '
$ws.Range("E29").Value = '  -0.10%  '
$ws.Range("D30").Value = '8.22'
$ws.Range("E30").Value = '  +3.83%  '
$ws.Range("D31").Value = '2.16'
$ws.Range("E31").Value = '  +4.59%  '
$ws.Range("D32").Value = '29.27'
$ws.Range("E32").Value = '  +9.29%  '
$ws.Range("D33").Value = '0.106'
$ws.Range("E34").Value = '  +1.44%  '
$ws.Range("E35").Value = '  +2.98%  '
$ws.Range("D36").Value = '6.01'
$ws.Range("E36").Value = '  +2.93%  '
$ws.Range("B37").Value = 'OKB'
$ws.Range("C37").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D37").Value = '53.30'
$ws.Range("E37").Value = '  +2.34%  '
$ws.Range("B38").Value = 'PEPE'
$ws.Range("C38").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D38").Value = '0.0₃0758'
$ws.Range("E38").Value = '  +7.04%  '
$ws.Range("D39").Value = '0.0400'
$ws.Range("E39").Value = '  +3.30%  '
$ws.Range("D40").Value = '433.64'
$ws.Range("E40").Value = '  +2.39%  '
$ws.Range("D41").Value = '3.029.88'
$ws.Range("E41").Value = '  +5.08%  '
$ws.Range("B42").Value = 'Cosmos'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D42").Value = '8.45'
$ws.Range("E42").Value = '  +2.22%  '
$ws.Range("B43").Value = 'dogwifhat'
$ws.Range("C43").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D43").Value = '2.78'
$ws.Range("E43").Value = '  +3.77%  '
$ws.Range("E44").Value = '  -5.27%  '
$ws.Range("E45").Value = '  +1.83%  '
$ws.Range("D46").Value = '2.22'
$ws.Range("E46").Value = '  +4.09%  '
$ws.Range("D47").Value = '26.53'
$ws.Range("E47").Value = '  +2.60%  '
$ws.Range("E48").Value = '  +0.11%  '
$ws.Range("D49").Value = '0.116'
$ws.Range("E49").Value = '  +2.33%  '
$ws.Range("E50").Value = '  +1.79%  '
$ws.Range("D51").Value = '35.32'
$ws.Range("E51").Value = '  +11.28%  '
